$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the date for the new timesheet entry row (18) - Dec 12, 2023 (serial 45272)
$ws.Range("A18").Value = 45272

# Enter the new formula for hours worked in column C, row 18
$ws.Range("C18").Formula = "=(1/60)*(19)"

# Move the active selection to C19, mirroring the user's next action after entry
$ws.Range("C19").Select()

$wb.Save()
